$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ScriptList")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1. Copy the old detail rows (2-7) from Sheet1 into Sheet2, starting at A1 (no header there)
$ws1.Range("A2:D7").Copy($ws2.Range("A1")) | Out-Null

# 2. On Sheet2 the "Run" column (D) should always read "N"
$ws2.Range("D1:D6").Value = "N"

# 3. Update Sheet1 row 2 with the new single remaining script entry
$ws1.Range("A2").Value = 7
$ws1.Range("B2").Value = "/src/test/resources/KW-Scripts/KWScript4.xlsx"
$ws1.Range("C2").Value = "Script1"
$ws1.Range("D2").Value = "Y"

# 4. Remove old rows 3-7 on Sheet1 so the used range shrinks back to A1:D2
$ws1.Range("A3:D7").Clear()

# 5. Update selections: Sheet2 first, then Sheet1 last so Sheet1 stays the active tab
$ws2.Range("A1:XFD6").Select() | Out-Null
$ws1.Range("C23").Select() | Out-Null
